$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.042.10"
$ws.Range("E2").Value = "  -0.85%  "

$ws.Range("D3").Value = "'1.904.70"
$ws.Range("E3").Value = "  -1.47%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'0.7409"
$ws.Range("E5").Value = "  -1.25%  "

$ws.Range("D6").Value = "'242.55"

$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "'0.3073"
$ws.Range("E8").Value = "  -3.50%  "

$ws.Range("D9").Value = "'26.23"
$ws.Range("E9").Value = "  -5.80%  "

$ws.Range("D10").Value = "'0.06902"
$ws.Range("E10").Value = "  -3.60%  "

$ws.Range("D11").Value = "'0.08069"
$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("D12").Value = "'0.7648"
$ws.Range("E12").Value = "  -2.21%  "

$ws.Range("D13").Value = "'1.902.00"
$ws.Range("E13").Value = "  -1.59%  "

$ws.Range("D14").Value = "'5.244"
$ws.Range("E14").Value = "  -2.83%  "

$ws.Range("D15").Value = "'91.51"
$ws.Range("E15").Value = "  -1.67%  "

$ws.Range("D16").Value = "'30.039.74"
$ws.Range("E16").Value = "  -0.83%  "

$ws.Range("D17").Value = "'6.081"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").Value = "'14.09"
$ws.Range("E18").Value = "  -3.23%  "

$ws.Range("D19").Value = "'0.000007774"
$ws.Range("E19").Value = "  -2.40%  "

$ws.Range("D20").Value = "'237.84"
$ws.Range("E20").Value = "  -5.51%  "

$ws.Range("D21").Value = "'2.173.88"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("D24").Value = "'7.089"
$ws.Range("E24").Value = "  +6.11%  "

$ws.Range("D25").Value = "'9.319"
$ws.Range("E25").Value = "  -2.47%  "

$ws.Range("D26").Value = "'166.47"
$ws.Range("E26").Value = "  +1.10%  "

$ws.Range("D27").Value = "'18.89"
$ws.Range("E27").Value = "  -1.10%  "

$ws.Range("D28").Value = "'0.1268"
$ws.Range("E28").Value = "  -2.42%  "

$ws.Range("D29").Value = "'2.040"
$ws.Range("E29").Value = "  -7.07%  "

$ws.Range("D30").Value = "'1.352"
$ws.Range("E30").Value = "  -1.71%  "

$ws.Range("E31").Value = "  -0.73%  "

$ws.Range("D32").Value = "'4.299"
$ws.Range("E32").Value = "  -2.72%  "

$ws.Range("D33").Value = "'4.058"
$ws.Range("E33").Value = "  -2.27%  "

$ws.Range("D34").Value = "'0.05427"
$ws.Range("E34").Value = "  +3.50%  "

$ws.Range("D35").Value = "'1.292"
$ws.Range("E35").Value = "  -2.53%  "

$ws.Range("D36").Value = "'0.7389"
$ws.Range("E36").Value = "  -2.49%  "

$ws.Range("E37").Value = "  -2.44%  "

$ws.Range("D38").Value = "'0.01957"
$ws.Range("E38").Value = "  +0.14%  "

$ws.Range("D39").Value = "'2.800"
$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("D40").Value = "'6.286"
$ws.Range("E40").Value = "  -3.47%  "

$ws.Range("D41").Value = "'0.4452"
$ws.Range("E41").Value = "  -1.69%  "

$ws.Range("D42").Value = "'73.71"
$ws.Range("E42").Value = "  -6.46%  "

$ws.Range("D43").Value = "'1.964"
$ws.Range("E43").Value = "  -0.78%  "

$ws.Range("D44").Value = "'1.003"
$ws.Range("E44").Value = "  +0.21%  "

$ws.Range("D45").Value = "'0.8356"
$ws.Range("E45").Value = "  -0.73%  "

$ws.Range("D46").Value = "'7.651"
$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("D47").Value = "'101.64"
$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("D48").Value = "'9.884"
$ws.Range("E48").Value = "  -1.12%  "

$ws.Range("D49").Value = "'2.064.01"
$ws.Range("E49").Value = "  -0.85%  "

$ws.Range("D50").Value = "'36.40"
$ws.Range("E50").Value = "  -3.61%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.1172"
$ws.Range("E51").Value = "  -3.68%  "
